# Commit: "Colocando header nos gráficos" (Adding header to the charts)
#
# For the first 4 sheets (same layout: Hidro, Gas Natural, Carvao, Nuclear,
# Oleos Comb, Biomassa, Eolica, Solar, Outros, Pot Compl, GD):
#   - Add a header label in A1 ("Fonte/Tecnologia") using the same style as
#     the existing B1:E1 header cells.
#   - Fix accentuation on several row labels (A2:A12).
#   - Remove the bold/centered header style from A2:A12 (now plain cells).
#
# Sheet 5 (Emissoes Totais): add A1 header "Período", fix accents on A2/A3,
# remove header style from A2/A3, and delete row 4 ("Teto") entirely.
#
# Sheet 6 (Custo Total): add A1 header "Tipo Expansão", change B1 from
# "Custo" to "2015", fix accents on A2/A3, remove header style from A2/A3,
# and update the B2/B3 values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share an identical layout.
# ---------------------------------------------------------------------
$rowLabels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

for ($s = 1; $s -le 4; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Add header cell A1, copying the style (including border/bold/center)
    # from the existing B1 header cell, then set its text.
    $ws.Range("B1").Copy($ws.Range("A1"))
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    foreach ($r in $rowLabels.Keys) {
        $cell = $ws.Cells.Item($r, 1)
        $cell.Value = $rowLabels[$r]
        $cell.ClearFormats()
    }
}

# ---------------------------------------------------------------------
# Sheet 5: Emissoes Totais (MtCO2eq)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy($ws5.Range("A1"))
$ws5.Range("A1").Value = "Período"

$ws5.Cells.Item(2, 1).Value = "P.Médio"
$ws5.Cells.Item(2, 1).ClearFormats()

$ws5.Cells.Item(3, 1).Value = "P.Crítico"
$ws5.Cells.Item(3, 1).ClearFormats()

# Remove row 4 ("Teto") entirely.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: Custo Total (bilhões de R$)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy($ws6.Range("A1"))
$ws6.Range("A1").Value = "Tipo Expansão"

# B1 changes from the text "Custo" to the text "2015". A plain
# Value assignment would have Excel auto-coerce the numeric-looking
# string into a real number, so force text via a quote-prefixed
# number format, then restore the original (un-prefixed) header
# style by pasting formats from another plain-text/style-1 cell.
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$ws6.Range("A2").Copy()
$ws6.Range("B1").PasteSpecial(-4122)

$ws6.Cells.Item(2, 1).Value = "Expansão Centralizada"
$ws6.Cells.Item(2, 1).ClearFormats()
$ws6.Cells.Item(2, 2).Value = 648

$ws6.Cells.Item(3, 1).Value = "Expansão por GD"
$ws6.Cells.Item(3, 1).ClearFormats()
$ws6.Cells.Item(3, 2).Value = 99
